$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.759.67"
$ws.Range("E2").Value = "  -0.02%  "

$ws.Range("D3").Value = "1.533.10"
$ws.Range("E3").Value = "  -1.97%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Formula = "'205.30"
$ws.Range("E5").Value = "  -0.43%  "

$ws.Range("E6").Value = "  -0.96%  "

$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Formula = "'0.244"
$ws.Range("E8").Value = "  -1.18%  "

$ws.Range("B9").Value = "Solana"
$ws.Range("C9").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D9").Formula = "'21.25"
$ws.Range("E9").Value = "  -3.21%  "

$ws.Range("E10").Value = "  -0.76%  "

$ws.Range("D11").Formula = "'0.0853"
$ws.Range("E11").Value = "  -0.84%  "

$ws.Range("D12").Value = "1.751.09"
$ws.Range("E12").Value = "  -1.96%  "

$ws.Range("D13").Value = "1.543.36"
$ws.Range("E13").Value = "  -1.32%  "

$ws.Range("E14").Value = "  -1.49%  "

$ws.Range("D15").Formula = "'0.505"
$ws.Range("E15").Value = "  -1.68%  "

$ws.Range("D16").Value = "26.741.27"

$ws.Range("D17").Formula = "'60.90"
$ws.Range("E17").Value = "  -0.99%  "

$ws.Range("D18").Formula = "'212.43"
$ws.Range("E18").Value = "  -0.69%  "

$ws.Range("D19").Value = "0.0₃0680"
$ws.Range("E19").Value = "  +0.85%  "

$ws.Range("E20").Value = "  -1.93%  "

$ws.Range("E21").Value = "  -0.06%  "

$ws.Range("E22").Value = "  -2.59%  "

$ws.Range("D23").Formula = "'9.07"
$ws.Range("E23").Value = "  -2.66%  "

$ws.Range("D24").Formula = "'1.94"
$ws.Range("E24").Value = "  -3.42%  "

$ws.Range("D25").Formula = "'151.33"
$ws.Range("E25").Value = "  -1.15%  "

$ws.Range("D26").Formula = "'6.56"
$ws.Range("E26").Value = "  -2.78%  "

$ws.Range("D27").Formula = "'14.76"
$ws.Range("E27").Value = "  -0.97%  "

$ws.Range("E28").Value = "  -0.09%  "

$ws.Range("E29").Value = "  -1.17%  "

$ws.Range("D30").Formula = "'1.10"
$ws.Range("E30").Value = "  -1.19%  "

$ws.Range("E31").Value = "  -2.21%  "

$ws.Range("E32").Value = "  +2.46%  "

$ws.Range("D33").Value = "1.361.18"
$ws.Range("E33").Value = "  -1.51%  "

$ws.Range("E34").Value = "  -0.15%  "

$ws.Range("E35").Value = "  -2.68%  "

$ws.Range("D36").Formula = "'0.947"
$ws.Range("E36").Value = "  +2.56%  "

$ws.Range("E37").Value = "  -0.30%  "

$ws.Range("E38").Value = "  -0.09%  "

$ws.Range("D39").Formula = "'0.519"
$ws.Range("E39").Value = "  -0.64%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Formula = "'5.72"
$ws.Range("E40").Value = "  +6.74%  "

$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Formula = "'0.799"
$ws.Range("E41").Value = "  -1.72%  "

$ws.Range("E42").Value = "  -0.06%  "

$ws.Range("D43").Formula = "'0.993"
$ws.Range("E43").Value = "  +0.09%  "

$ws.Range("E44").Value = "  +0.52%  "

$ws.Range("D45").Formula = "'1.74"
$ws.Range("E45").Value = "  -2.11%  "

$ws.Range("D46").Formula = "'62.48"
$ws.Range("E46").Value = "  -1.14%  "

$ws.Range("D47").Value = "1.665.30"
$ws.Range("E47").Value = "  -2.03%  "

$ws.Range("D48").Formula = "'85.11"
$ws.Range("E48").Value = "  -0.35%  "

$ws.Range("E49").Value = "  +2.19%  "

$ws.Range("D50").Value = "0.0₇0970"
$ws.Range("E50").Value = "  -1.64%  "

$ws.Range("D51").Formula = "'0.0941"
$ws.Range("E51").Value = "  -0.85%  "
